# Prix Spot sheet: a new daily price column ("15-jun") is inserted before
# the existing "16-jun" column, shifting the old "16-jun"/"17-jun" data one
# column to the right (B->C, C->D). The new "15-jun" column is then filled
# with its own price series.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Prix Spot")

# Insert a new column at B; this shifts the old B ("16-jun") and C
# ("17-jun") columns to C and D, carrying their data/formatting along.
$ws.Columns.Item(2).Insert()

# The insert copies column A's header/bold/border style onto every new
# B cell; strip that back off the data rows (row 1 keeps the header style).
$ws.Range("B2:B25").ClearFormats()

# New column header.
$ws.Range("B1").Value = "15-jun"

# New "15-jun" price series (hour rows 2-25 <-> 00-01 .. 23-24).
$values = @(
    51.35,
    28.31,
    26.87,
    21.88,
    18.78,
    17.86,
    18.47,
    15.13,
    8.970000000000001,
    4.55,
    0,
    -0.02,
    -1.21,
    -5.6,
    -5,
    -2,
    -0.01,
    -0.01,
    12.37,
    19.29,
    33.96,
    39.96,
    61.7,
    53.03
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $ws.Cells.Item(2 + $i, 2).Value = $values[$i]
}
